$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.999.48'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '1.817.58'
$ws.Range('E3').Value = '  +2.45%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '337.22'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9996'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  +11.43%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3521'
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '45.59'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.148'
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07458'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '23.01'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.003'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.275'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.329'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').Value = '1.811.66'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001087'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06706'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '82.28'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.28'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.419'
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').Value = '28.055.71'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.89'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.409'
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.472'
$ws.Range('E26').Value = '  +3.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.76'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '155.33'
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').Value = '2.022.13'
$ws.Range('E29').Value = '  +2.18%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.305'
$ws.Range('E30').Value = '  -8.84%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '132.78'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.074'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.983'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.09155'
$ws.Range('E34').Value = '  +2.67%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.40'
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02373'
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6719'
$ws.Range('E37').Value = '  -1.57%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06296'
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.246'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2176'
$ws.Range('E40').Value = '  +0.91%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.501'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.219'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.090'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '14.28'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6152'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.874'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '128.21'
$ws.Range('E48').Value = '  -3.19%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.055'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.183'
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07136'
$ws.Range('E51').Value = '  -4.78%  '
